$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell's text (Price / Volume(1h) columns, plus the
# Stacks/Aave row swap in rows 49-50). NumberFormat is forced to "@" (Text)
# before the write so numeric-looking strings (e.g. "1.00", "0.111") are not
# silently coerced to numbers by Excel's type inference, then ClearFormats()
# removes that temporary formatting again so the cell keeps its original
# (unstyled) look - matching how the sheet was authored originally.
$changes = @(
    @{ Cell = 'D2'; Value = '40.046.81' }
    @{ Cell = 'E2'; Value = '  +0.03%  ' }
    @{ Cell = 'D3'; Value = '2.221.51' }
    @{ Cell = 'E3'; Value = '  -0.73%  ' }
    @{ Cell = 'E4'; Value = '  +0.12%  ' }
    @{ Cell = 'D5'; Value = '289.87' }
    @{ Cell = 'E5'; Value = '  -0.80%  ' }
    @{ Cell = 'D6'; Value = '88.16' }
    @{ Cell = 'E6'; Value = '  +2.76%  ' }
    @{ Cell = 'E7'; Value = '  -0.82%  ' }
    @{ Cell = 'E8'; Value = '  +0.04%  ' }
    @{ Cell = 'E9'; Value = '  +0.50%  ' }
    @{ Cell = 'D10'; Value = '30.52' }
    @{ Cell = 'E10'; Value = '  +0.04%  ' }
    @{ Cell = 'E11'; Value = '  -2.35%  ' }
    @{ Cell = 'D12'; Value = '0.111' }
    @{ Cell = 'E12'; Value = '  +3.20%  ' }
    @{ Cell = 'D13'; Value = '6.50' }
    @{ Cell = 'E13'; Value = '  +1.83%  ' }
    @{ Cell = 'D14'; Value = '2.566.21' }
    @{ Cell = 'E14'; Value = '  -0.64%  ' }
    @{ Cell = 'E15'; Value = '  -1.60%  ' }
    @{ Cell = 'D16'; Value = '2.218.77' }
    @{ Cell = 'E16'; Value = '  -0.43%  ' }
    @{ Cell = 'E17'; Value = '  +0.72%  ' }
    @{ Cell = 'D18'; Value = '39.998.93' }
    @{ Cell = 'E18'; Value = '  +0.23%  ' }
    @{ Cell = 'D19'; Value = '11.61' }
    @{ Cell = 'E19'; Value = '  +8.83%  ' }
    @{ Cell = 'E20'; Value = '  -0.83%  ' }
    @{ Cell = 'D21'; Value = '5.82' }
    @{ Cell = 'E21'; Value = '  +0.29%  ' }
    @{ Cell = 'D22'; Value = '65.69' }
    @{ Cell = 'E22'; Value = '  +0.29%  ' }
    @{ Cell = 'D23'; Value = '235.88' }
    @{ Cell = 'E23'; Value = '  +1.10%  ' }
    @{ Cell = 'D24'; Value = '1.00' }
    @{ Cell = 'E24'; Value = '  +0.00%  ' }
    @{ Cell = 'D25'; Value = '2.46' }
    @{ Cell = 'E25'; Value = '  +1.86%  ' }
    @{ Cell = 'E26'; Value = '  -1.31%  ' }
    @{ Cell = 'D27'; Value = '22.61' }
    @{ Cell = 'E27'; Value = '  -1.98%  ' }
    @{ Cell = 'E28'; Value = '  -0.50%  ' }
    @{ Cell = 'D29'; Value = '9.22' }
    @{ Cell = 'E29'; Value = '  -0.49%  ' }
    @{ Cell = 'D30'; Value = '155.24' }
    @{ Cell = 'E30'; Value = '  +0.01%  ' }
    @{ Cell = 'D31'; Value = '31.85' }
    @{ Cell = 'E31'; Value = '  -4.60%  ' }
    @{ Cell = 'E33'; Value = '  +2.03%  ' }
    @{ Cell = 'D34'; Value = '0.0718' }
    @{ Cell = 'E34'; Value = '  +1.09%  ' }
    @{ Cell = 'E35'; Value = '  +0.28%  ' }
    @{ Cell = 'D36'; Value = '2.87' }
    @{ Cell = 'E36'; Value = '  +6.69%  ' }
    @{ Cell = 'E37'; Value = '  -0.47%  ' }
    @{ Cell = 'D38'; Value = '15.81' }
    @{ Cell = 'E38'; Value = '  -5.01%  ' }
    @{ Cell = 'D39'; Value = '0.0986' }
    @{ Cell = 'E39'; Value = '  +0.25%  ' }
    @{ Cell = 'D40'; Value = '1.70' }
    @{ Cell = 'E40'; Value = '  +1.98%  ' }
    @{ Cell = 'D41'; Value = '2.104.61' }
    @{ Cell = 'E41'; Value = '  +7.54%  ' }
    @{ Cell = 'D42'; Value = '3.83' }
    @{ Cell = 'E42'; Value = '  +1.88%  ' }
    @{ Cell = 'D43'; Value = '2.14' }
    @{ Cell = 'E43'; Value = '  -2.18%  ' }
    @{ Cell = 'E44'; Value = '  -0.77%  ' }
    @{ Cell = 'D45'; Value = '9.89' }
    @{ Cell = 'E45'; Value = '  +4.01%  ' }
    @{ Cell = 'D46'; Value = '17.53' }
    @{ Cell = 'E46'; Value = '  +6.89%  ' }
    @{ Cell = 'E47'; Value = '  +2.16%  ' }
    @{ Cell = 'D48'; Value = '2.432.98' }
    @{ Cell = 'E48'; Value = '  -0.61%  ' }
    @{ Cell = 'B49'; Value = 'Aave' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D49'; Value = '88.95' }
    @{ Cell = 'E49'; Value = '  -0.13%  ' }
    @{ Cell = 'B50'; Value = 'Stacks' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D50'; Value = '1.44' }
    @{ Cell = 'E50'; Value = '  -1.01%  ' }
    @{ Cell = 'D51'; Value = '69.04' }
    @{ Cell = 'E51'; Value = '  -2.96%  ' }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
    $cell.ClearFormats()
}
